# B6-PowerPoint.pptx edit:
#  1. Swap the table style applied to the three "balance sheet" tables
#     (slides 14, 15, 16) from {49E1687B-...} to {E7B7E081-...}.
#  2. Swap the (single, shared) master theme's colour scheme from the
#     "Integral / Red Violet" palette to the "Office Theme / Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style swap -----------------------------------------------
$oldStyle = "{49E1687B-9BE0-435F-812B-205CE9CC7F9B}"
$newStyle = "{E7B7E081-1ADB-4F54-8B0E-60081873C8DF}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyle, $true)
        }
    }
}

# --- 2. Theme colour scheme swap (Integral -> Office Theme) -----------
$scheme = $p.SlideMaster.Theme.ThemeColorScheme
# Index: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
#        8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# Values are packed BGR integers (VBA RGB() order) for the target
# "Office" colour scheme.
$scheme.Colors(1).RGB  = 0          # dk1     000000
$scheme.Colors(2).RGB  = 16777215   # lt1     FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2     44546A
$scheme.Colors(4).RGB  = 15132391   # lt2     E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1 5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2 ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3 A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4 FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5 4472C4
$scheme.Colors(10).RGB = 4697456    # accent6 70AD47
$scheme.Colors(11).RGB = 12673797   # hlink   0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
